# Update countries & provincias Spain
# -----------------------------------------------------------------------
# This script applies the 28-Abril-2020 21:52 data refresh to the "Pais"
# sheet:
#   1. Updates the per-country case figures that changed since 21:22.
#   2. Refreshes Guinea Ecuatorial's and Maldivas' rows with their new
#      figures (Guinea Ecuatorial's total overtakes Isla de Man /
#      Tanzania / Vietnam in the ranking).
#   3. Re-sorts the country table (A4:H216) by "Casos totales" (column B)
#      descending, same as the source sheet is always kept in.
#   4. Updates the "Datos actualizados a ..." footer timestamp.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-CountryRow($name) {
    # LookAt:=xlWhole (1) avoids ambiguity between e.g. "Guinea" / "Guinea Ecuatorial" / "Guinea-Bisau"
    $found = $ws.Range("A4:A216").Find($name, [System.Type]::Missing, [System.Type]::Missing, 1)
    return $found.Row
}

function Set-CountryStats($name, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $r = Get-CountryRow $name
    if ($total -ne $null)        { $ws.Cells.Item($r, 2).Value = $total }
    if ($nuevos -ne $null)       { $ws.Cells.Item($r, 3).Value = $nuevos }
    if ($activos -ne $null)      { $ws.Cells.Item($r, 4).Value = $activos }
    if ($recuperados -ne $null)  { $ws.Cells.Item($r, 5).Value = $recuperados }
    if ($criticos -ne $null)     { $ws.Cells.Item($r, 6).Value = $criticos }
    if ($muertesHoy -ne $null)   { $ws.Cells.Item($r, 7).Value = $muertesHoy }
    if ($muertes -ne $null)      { $ws.Cells.Item($r, 8).Value = $muertes }
}

# --- Per-country figure refresh ----------------------------------------
Set-CountryStats "Estados Unidos" 1026771 16415 140138 828364 $null 1472 58269
Set-CountryStats "Francia"        $null   $null $null  96669  $null 367  23660
Set-CountryStats "Alemania"       159431  673   $null  35816  $null 89   6215
Set-CountryStats "Brasil"         68289   1788  $null  32464  $null 140  4683
Set-CountryStats "Canada"         49815   1315  19047  27916  $null 145  2852
Set-CountryStats "India"          31360   1909  7747   22605  $null 69   1008
Set-CountryStats "Suiza"          $null   $null 22600  4965   $null $null $null
Set-CountryStats "Israel"         15728   173   7746   7772   116   6    210
Set-CountryStats "Sudafrica"      4996    203   2073   2830   $null 3    93

# Guinea Ecuatorial & Maldivas figures
Set-CountryStats "Guinea Ecuatorial" 315 57 9  305 0    0    1
Set-CountryStats "Maldivas"          250 24 $null 233 $null $null $null

# --- Re-sort the country table by total cases, descending --------------
$sortRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$sortRange.Sort($sortKey, 2)

# --- Footer timestamp ----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 21:52"
